$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4, 16 and 28 in the original data are exact duplicates of the row
# immediately above them (TURNER EQUITY INVESTORS INC, NASTA INTERNATIONAL
# INC and TRI-LITE INC respectively). Remove those duplicate rows; Excel
# shifts everything below each deleted row up by one, and the dimension
# shrinks accordingly (from A1:G123 to A1:G120).
# Delete from the bottom up so earlier row numbers stay valid.
$ws.Rows.Item(28).Delete()
$ws.Rows.Item(16).Delete()
$ws.Rows.Item(4).Delete()
